$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the UBID values in column E (rows 2-15)
$ws.Range("E2").Value = "86HJQCC9+5JJ-2-3-2-3"
$ws.Range("E3").Value = "86HJX5QV+FJ3-2-3-2-2"
$ws.Range("E4").Value = "86HJQ8Q5+R6V-1-2-1-1"
$ws.Range("E5").Value = "86HJX6JP+H99-1-1-1-2"
$ws.Range("E6").Value = "86HJQ76M+883-1-2-1-1"
$ws.Range("E7").Value = "86HJW825+V3M-2-2-3-1"
$ws.Range("E8").Value = "86HJX6GX+F4G-2-4-2-3"
$ws.Range("E9").Value = "86HJX66G+P7C-2-3-2-3"
$ws.Range("E10").Value = "86HJM8JW+XMV-1-4-1-3"
$ws.Range("E11").Value = "86HJPCWJ+R59-1-5-2-4"
$ws.Range("E12").Value = "86HJR7QR+98F-2-1-1-1"
$ws.Range("E13").Value = "86HJQ9R3+FHW-1-2-0-3"
$ws.Range("E14").Value = "86HJW5RW+VGV-1-2-2-2"
$ws.Range("E15").Value = "86HJX838+8M7-1-3-1-2"

# Materialize a few new, empty rows below the data (rows 36-38), matching
# a user having touched those cells (e.g. while reviewing/extending the
# sheet) without putting any value into them.
$ws.Range("E36").Interior.Pattern = -4142
$ws.Range("E37").Interior.Pattern = -4142
$ws.Range("E38").Interior.Pattern = -4142

# Update the active selection
$ws.Range("E11").Select() | Out-Null
